$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 122644.96
$ws.Range("I17").Value = 80.5
$ws.Range("J17").Value = 127860.47
$ws.Range("K17").Value = 241.5
$ws.Range("L17").Value = 383581.41
$ws.Range("M17").Value = -73.5
$ws.Range("N17").Value = -383917.41

$ws.Range("H99").Value = 1328
$ws.Range("J99").Value = 3700
$ws.Range("L99").Value = 11100
$ws.Range("N99").Value = -14096

$ws.Range("H112").Value = 1111.1111
$ws.Range("J112").Value = 1212.5
$ws.Range("L112").Value = 3637.5
$ws.Range("N112").Value = -5853.5

$ws.Range("H127").Value = 66667344
$ws.Range("I127").Value = 83333830
$ws.Range("K127").Value = 250001490
$ws.Range("M127").Value = -249996530

$ws.Range("H129").Value = 939.4386
$ws.Range("J129").Value = 1029.3334
$ws.Range("L129").Value = 3088.0002
$ws.Range("N129").Value = -13088.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2798.9324
$ws.Range("I32").Value = 2508.5933
$ws.Range("K32").Value = 2508.5933
$ws.Range("M32").Value = -2221.5933

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 99342
$ws.Range("J47").Value = 99342
$ws.Range("L47").Value = 99342
$ws.Range("N47").Value = -100382

$ws.Range("H137").Value = 43285.668
$ws.Range("J137").Value = 43285.668
$ws.Range("L137").Value = 43285.668
$ws.Range("N137").Value = -53485.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 248.75
$ws.Range("I22").Value = 298.33334
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 298.33334
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 51.66665999999998
$ws.Range("N22").Value = -800

$ws.Range("H62").Value = 2939.2856
$ws.Range("I62").Value = 2312.5
$ws.Range("J62").Value = 3775
$ws.Range("K62").Value = 2312.5
$ws.Range("L62").Value = 3775
$ws.Range("M62").Value = -1688.5
$ws.Range("N62").Value = -5023

$ws.Range("H65").Value = 2939.2856
$ws.Range("I65").Value = 2312.5
$ws.Range("J65").Value = 3775
$ws.Range("K65").Value = 11562.5
$ws.Range("L65").Value = 18875
$ws.Range("M65").Value = -8442.5
$ws.Range("N65").Value = -25115

$ws.Range("H93").Value = 10067.667
$ws.Range("I93").Value = 6081.2
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 6081.2
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -4209.2
$ws.Range("N93").Value = -33744

$ws.Range("H95").Value = 15584.6
$ws.Range("J95").Value = 15584.6
$ws.Range("L95").Value = 15584.6
$ws.Range("N95").Value = -21076.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2854.5715
$ws.Range("I19").Value = 1500
$ws.Range("J19").Value = 4660.6665
$ws.Range("K19").Value = 4500
$ws.Range("L19").Value = 13981.9995
$ws.Range("M19").Value = -4326
$ws.Range("N19").Value = -14329.9995

$ws.Range("H70").Value = 5675
$ws.Range("I70").Value = 4208.6665
$ws.Range("J70").Value = 6774.75
$ws.Range("K70").Value = 12625.9995
$ws.Range("L70").Value = 20324.25
$ws.Range("M70").Value = -12310.9995
$ws.Range("N70").Value = -20954.25

$ws.Range("H73").Value = 5675
$ws.Range("I73").Value = 4208.6665
$ws.Range("J73").Value = 6774.75
$ws.Range("K73").Value = 12625.9995
$ws.Range("L73").Value = 20324.25
$ws.Range("M73").Value = -11533.9995
$ws.Range("N73").Value = -22508.25

$ws.Range("H109").Value = 2056.8823
$ws.Range("J109").Value = 3124.5
$ws.Range("L109").Value = 9373.5
$ws.Range("N109").Value = -11453.5

$ws.Range("H113").Value = 574.871
$ws.Range("J113").Value = 527.04346
$ws.Range("L113").Value = 1581.13038
$ws.Range("N113").Value = -5921.130380000001

$ws.Range("H121").Value = 1023.96295
$ws.Range("I121").Value = 426.5
$ws.Range("J121").Value = 1194.6666
$ws.Range("K121").Value = 1279.5
$ws.Range("L121").Value = 3583.9998
$ws.Range("M121").Value = 30.5
$ws.Range("N121").Value = -6203.9998

$ws.Range("H131").Value = 870.1340300000001
$ws.Range("J131").Value = 901.2111
$ws.Range("L131").Value = 2703.6333
$ws.Range("N131").Value = -12783.6333

$ws.Range("H134").Value = 5592.6787
$ws.Range("I134").Value = 2872.4546
$ws.Range("J134").Value = 7352.8237
$ws.Range("K134").Value = 8617.363799999999
$ws.Range("L134").Value = 22058.4711
$ws.Range("M134").Value = -3547.363799999999
$ws.Range("N134").Value = -32198.4711

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1339.5862
$ws.Range("I102").Value = 1276.8
$ws.Range("J102").Value = 1479.1111
$ws.Range("K102").Value = 1276.8
$ws.Range("L102").Value = 1479.1111
$ws.Range("M102").Value = 345.2
$ws.Range("N102").Value = -4723.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3172.182
$ws.Range("I40").Value = 3172.182
$ws.Range("K40").Value = 3172.182
$ws.Range("M40").Value = -3036.182

$ws.Range("H93").Value = 34673.094
$ws.Range("I93").Value = 1484
$ws.Range("J93").Value = 175726.75
$ws.Range("K93").Value = 1484
$ws.Range("L93").Value = 175726.75
$ws.Range("M93").Value = -236
$ws.Range("N93").Value = -178222.75

$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws.Range("H101").Value = 23697.2
$ws.Range("J101").Value = 23697.2
$ws.Range("L101").Value = 23697.2
$ws.Range("N101").Value = -30187.2

$ws.Range("H122").Value = 3532.6428
$ws.Range("I122").Value = 3080
$ws.Range("J122").Value = 3872.125
$ws.Range("K122").Value = 9240
$ws.Range("L122").Value = 11616.375
$ws.Range("M122").Value = -6790
$ws.Range("N122").Value = -16516.375

$ws.Range("H136").Value = 2166.1333
$ws.Range("I136").Value = 2181.0908
$ws.Range("J136").Value = 2125
$ws.Range("K136").Value = 6543.2724
$ws.Range("L136").Value = 6375
$ws.Range("M136").Value = -3993.2724
$ws.Range("N136").Value = -11475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1794.8235
$ws.Range("I126").Value = 1251
$ws.Range("J126").Value = 1962.1538
$ws.Range("K126").Value = 3753
$ws.Range("L126").Value = 5886.4614
$ws.Range("M126").Value = -1283
$ws.Range("N126").Value = -10826.4614
